# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Periodo Mora" (column E) list is reversed from newest-first
# (2009 .. 1608) to oldest-first (1608 .. 2009), and the "Valor Mora"
# (column F) figures are re-pointed to match the re-ordered periods.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$rows = @(
    [PSCustomObject]@{ Row = 16; E = "1608"; F = 27556 },
    [PSCustomObject]@{ Row = 17; E = "1609"; F = 27556 },
    [PSCustomObject]@{ Row = 18; E = "1610"; F = 27556 },
    [PSCustomObject]@{ Row = 19; E = "1611"; F = 27556 },
    [PSCustomObject]@{ Row = 20; E = "1612"; F = 27556 },
    [PSCustomObject]@{ Row = 21; E = "1701"; F = 27556 },
    [PSCustomObject]@{ Row = 22; E = "1702"; F = 27556 },
    [PSCustomObject]@{ Row = 23; E = "1704"; F = 27556 },
    [PSCustomObject]@{ Row = 24; E = "1705"; F = 27556 },
    [PSCustomObject]@{ Row = 25; E = "1706"; F = 27556 },
    [PSCustomObject]@{ Row = 26; E = "1707"; F = 27556 },
    [PSCustomObject]@{ Row = 27; E = "1708"; F = 27556 },
    [PSCustomObject]@{ Row = 28; E = "1709"; F = 27556 },
    [PSCustomObject]@{ Row = 29; E = "1710"; F = 27556 },
    [PSCustomObject]@{ Row = 30; E = "1711"; F = 27556 },
    [PSCustomObject]@{ Row = 31; E = "1712"; F = 27556 },
    [PSCustomObject]@{ Row = 32; E = "1801"; F = 27556 },
    [PSCustomObject]@{ Row = 33; E = "1802"; F = 27556 },
    [PSCustomObject]@{ Row = 34; E = "1803"; F = 27556 },
    [PSCustomObject]@{ Row = 35; E = "1804"; F = 27556 },
    [PSCustomObject]@{ Row = 36; E = "1805"; F = 27556 },
    [PSCustomObject]@{ Row = 37; E = "1806"; F = 27556 },
    [PSCustomObject]@{ Row = 38; E = "1807"; F = 27556 },
    [PSCustomObject]@{ Row = 39; E = "1808"; F = 27556 },
    [PSCustomObject]@{ Row = 40; E = "1809"; F = 31249 },
    [PSCustomObject]@{ Row = 41; E = "1810"; F = 31249 },
    [PSCustomObject]@{ Row = 42; E = "1811"; F = 31249 },
    [PSCustomObject]@{ Row = 43; E = "1812"; F = 31249 },
    [PSCustomObject]@{ Row = 44; E = "1901"; F = 31249 },
    [PSCustomObject]@{ Row = 45; E = "1902"; F = 31249 },
    [PSCustomObject]@{ Row = 46; E = "1903"; F = 31249 },
    [PSCustomObject]@{ Row = 47; E = "1904"; F = 31249 },
    [PSCustomObject]@{ Row = 48; E = "1905"; F = 31249 },
    [PSCustomObject]@{ Row = 49; E = "1906"; F = 31249 },
    [PSCustomObject]@{ Row = 50; E = "1907"; F = 31249 },
    [PSCustomObject]@{ Row = 51; E = "1908"; F = 31249 },
    [PSCustomObject]@{ Row = 52; E = "1909"; F = 31249 },
    [PSCustomObject]@{ Row = 53; E = "1910"; F = 31249 },
    [PSCustomObject]@{ Row = 54; E = "1911"; F = 31249 },
    [PSCustomObject]@{ Row = 55; E = "1912"; F = 31249 },
    [PSCustomObject]@{ Row = 56; E = "2001"; F = 31249 },
    [PSCustomObject]@{ Row = 57; E = "2002"; F = 31249 },
    [PSCustomObject]@{ Row = 58; E = "2003"; F = 31249 },
    [PSCustomObject]@{ Row = 59; E = "2004"; F = 31249 },
    [PSCustomObject]@{ Row = 60; E = "2005"; F = 31249 },
    [PSCustomObject]@{ Row = 61; E = "2006"; F = 31249 },
    [PSCustomObject]@{ Row = 62; E = "2007"; F = 31249 },
    [PSCustomObject]@{ Row = 63; E = "2008"; F = 31249 },
    [PSCustomObject]@{ Row = 64; E = "2009"; F = 29166 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}
